# Redesign the product report: new title/header rows, renamed sheet,
# a new "Condones" product row, and a blank styled footer row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab ---------------------------------------------
$ws.Name = "Reporte Productos"

# --- Snapshot the existing table (rows 1-4) before we overwrite it --------
$oldRow2 = @($ws.Range("A2").Value2, $ws.Range("B2").Value2, $ws.Range("C2").Value2, $ws.Range("D2").Value2, $ws.Range("E2").Value2, $ws.Range("F2").Value2, $ws.Range("G2").Value2)
$oldRow3 = @($ws.Range("A3").Value2, $ws.Range("B3").Value2, $ws.Range("C3").Value2, $ws.Range("D3").Value2, $ws.Range("E3").Value2, $ws.Range("F3").Value2, $ws.Range("G3").Value2)
$oldRow4 = @($ws.Range("A4").Value2, $ws.Range("B4").Value2, $ws.Range("C4").Value2, $ws.Range("D4").Value2, $ws.Range("E4").Value2, $ws.Range("F4").Value2, $ws.Range("G4").Value2)

# --- Clear the whole former table so we can rebuild it from scratch -------
$ws.Range("A1:G4").ClearContents()

# --- Row 1: report title + date banner (merged, bold 12pt) ----------------
$ws.Range("A1").Value = "Reporte de Ventas"
$ws.Range("D1").Value = "Fecha  2021-11-18 15:40:09"
$ws.Range("A1:C1").Merge()
$ws.Range("D1:F1").Merge()
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").Font.Size = 12

# Row 2 is intentionally left blank (spacer row).

# --- Row 3: column sub-headers (merged, bold 12pt) -------------------------
$ws.Range("A3").Value = "Nombre"
$ws.Range("B3").Value = "Concentración"
$ws.Range("B3:C3").Merge()
$ws.Range("A3:B3").Font.Bold = $true
$ws.Range("A3:B3").Font.Size = 12

# --- Row 4: first product (renamed from "buscapina ll") -------------------
$ws.Range("A4").Value = "buscapina ejemplo"
$ws.Range("B4").Value = $oldRow2[1]
$ws.Range("C4").Value = $oldRow2[2]
$ws.Range("D4").Value = $oldRow2[3]
$ws.Range("E4").Value = $oldRow2[4]
$ws.Range("F4").Value = $oldRow2[5]
$ws.Range("G4").Value = $oldRow2[6]

# --- Row 5: brand-new "Condones" product -----------------------------------
$ws.Range("A5").Value = "Condones"
$ws.Range("B5").Value = "condones "
$ws.Range("C5").Value = "c"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "asdfg"
$ws.Range("F5").Value = "inyeccion"
$ws.Range("G5").Value = "local host"

# --- Row 6: second original product, unchanged ------------------------------
$ws.Range("A6").Value = $oldRow3[0]
$ws.Range("B6").Value = $oldRow3[1]
$ws.Range("C6").Value = $oldRow3[2]
$ws.Range("D6").Value = $oldRow3[3]
$ws.Range("E6").Value = $oldRow3[4]
$ws.Range("F6").Value = $oldRow3[5]
$ws.Range("G6").Value = $oldRow3[6]

# --- Row 7: third original product, unchanged --------------------------------
$ws.Range("A7").Value = $oldRow4[0]
$ws.Range("B7").Value = $oldRow4[1]
$ws.Range("C7").Value = $oldRow4[2]
$ws.Range("D7").Value = $oldRow4[3]
$ws.Range("E7").Value = $oldRow4[4]
$ws.Range("F7").Value = $oldRow4[5]
$ws.Range("G7").Value = $oldRow4[6]

# --- Row 8: blank footer row, centered ---------------------------------------
$ws.Range("A8:G8").HorizontalAlignment = -4108

$ws.Range("G8").Select()
